$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Test Yearly").Name = "Test_Yearly"
$wb.Worksheets.Item("Test Weekly").Name = "Test_Weekly"
